$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.608888149261475
$ws.Range("B1").Value = 4.235321044921875
$ws.Range("C1").Value = 3.624852657318115
$ws.Range("D1").Value = 1.805981636047363
$ws.Range("E1").Value = 1.04540228843689
